# Adds a new "FE achievement aims" metric row to the dataText sheet.
# This inserts a new row 15 (pushing the existing rows 15-24 down to 16-25)
# and populates it with the new metric's text, mirroring the layout/style
# of the existing "FE achievements" row (row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 15, shifting everything below it down ---
$ws.Rows("15:15").Insert()

# Copy the formatting (styles) from row 14 (an existing "FE achievements"
# row with the same column layout) into the new row 15 so the new row
# matches the workbook's existing look (border, alignment, fonts, etc.)
$ws.Range("A14:M14").Copy()
$ws.Range("A15:M15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row height for the new row
$ws.Rows("15:15").RowHeight = 348.5

# --- Populate the new row's content ---
$ws.Range("A15").Value2 = "achievementAims"
$ws.Range("B15").Value2 = "AY22/23 data"
$ws.Range("C15").Value2 = " Further education and skills include 19+ apprenticeships and publicly-funded adult (19+) learning, including community learning, delivered by an FE institution, a training provider or within a local community."
$ws.Range("D15").Value2 = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/further-education-and-skills/2022-23'>Individualised Learner Record</a>"

$ws.Range("E15").Value2 = @'
Further education and skills included are 19+ apprenticeships and publicly-funded adult learning, including community learning, delivered by an FE institution, a training provider or within a local community.
FE and skills does not includer higher education, unless delivered as part of an apprenticeship programme.
Apprenticeships are paid jobs that incorporate on-the-job and off-the-job training leading to nationally recognised qualifications.
Community learning funds a wide range of non-formal courses (e.g. IT or employability skills) and activity targeted at deprived areas or disadvantaged groups. They can be offered by local authorities, colleges, community groups.
Achievements are the number of programme aims successfully completed in an invidual aim in an academic year.
'@

$ws.Range("F15").Value2 = @'
<ol>
  <li>Total achievements is the count of programme level aims completed at any point during the stated academic period.</li>
<li> Education and Training and Apprenticeship aims are 19 plus. </li>
 <li>Years shown represent academic years.</li>
<li>Use caution when interpreting this data. A difference between subgroups does not necessarily imply any causality. There could be other contributing factors at work.</li>
</ol>
'@

$ws.Range("G15").Value2 = "FE achievements aims"
$ws.Range("H15").Value2 = "are FE achievement aim volumes changing"
$ws.Range("I15").Value2 = "The number of FE achievement aims in"
$ws.Range("J15").Value2 = "FE achievement aims"
$ws.Range("K15").Value2 = "FE achievement aims"
$ws.Range("L15").Value2 = "share of FE achievement aims"
$ws.Range("M15").Value2 = "FE achievement aims are"

# --- Update the sheet view to match where the author left the selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C16").Select()
